$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27; this shifts the existing rows 27-40 down
# to become rows 28-41, preserving their data and formatting.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new weekly data entry.
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "Macroferia Regional de Talca"
$ws.Range("C27").Value = "Maule"
$ws.Range("D27").Value = 44438
$ws.Range("D27").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = 100112013
$ws.Range("G27").Value = "Alcachofa"
$ws.Range("H27").Value = "Madrigal"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("N27").Value = "`$/caja 40 unidades"
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("P27").Value = 250
$ws.Range("Q27").Value = 40
$ws.Range("R27").Value = "Hortaliza"
